# Apply the "Initial commit for kfserving abtesting" edit to the single
# slide of the request-routing deck.
#
# Summary of changes on Slide 1:
#   1. "Rectangle 10"  (background panel) moves down / shrinks:
#        Top    1086522 EMU -> 1389528 EMU
#        Height 5077610 EMU -> 4518745 EMU
#   2. "TextBox 4" ("Istio Virtual Service" label) is repositioned/resized,
#      re-wrapped, centred, and re-worded to
#      "Ingress / custom traffic resource".
#   3. Shape "Rounded Rectangle 90" (an empty decorative rounded rect) is
#      removed entirely.
#   4. Shapes "TextBox 142" ("Objectives" label) and "Graphic 143" (the
#      checkbox icon) are removed entirely.
#   5. The slide's Footer and Slide Number placeholders are turned off,
#      which removes their placeholder shapes from the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Resize/move the big background rectangle -------------------------
$rect = $s.Shapes.Item("Rectangle 10")
$rect.Top    = 109.41165924072266   # -> 1389528 EMU
$rect.Height = 355.80670166015625   # -> 4518745 EMU

# --- 2. Move/resize/re-wrap/re-centre/re-word the label textbox ----------
$lbl = $s.Shapes.Item("TextBox 4")
$lbl.Left   = 138.28323364257812    # -> 1756197 EMU
$lbl.Top    = 294.7419128417969     # -> 3743222 EMU
$lbl.Width  = 136.28402709960938    # -> 1730807 EMU
$lbl.Height = 50.892208099365234    # -> 646331 EMU

$lbl.TextFrame.WordWrap = -1
$lbl.TextFrame.TextRange.Text = "Ingress / custom traffic resource"
$lbl.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- 3/4. Delete the now-unused decorative shapes -------------------------
$s.Shapes.Item("Rounded Rectangle 90").Delete()
$s.Shapes.Item("TextBox 142").Delete()
$s.Shapes.Item("Graphic 143").Delete()

# --- 5. Remove the footer / slide-number placeholders ----------------------
# A plain Shape.Delete() on these special placeholders only clears their
# text (PowerPoint re-creates an empty placeholder on save); toggling the
# slide's HeadersFooters visibility actually drops the shapes.
$hf = $s.HeadersFooters
$hf.Footer.Visible = 0
$hf.SlideNumber.Visible = 0
